$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1, copying the format of the existing header cell (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Column I values (rows 2-38)
$iValues = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,6,7,1,1)
# Column J values (rows 2-38)
$jValues = @(6,4,6,6,7,3,5,3,6,6,7,7,4,8,6,7,8,5,6,7,6,6,5,5,7,6,8,6,6,6,6,5,6,9,9,3,2)

for ($r = 2; $r -le 38; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
